$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in notes cell (G162): "repreated" -> "repeated"
$ws.Range("G162").Value = "For plots deleted points with 6 highest alphas or 6 highest lambdas (if repeated)"

# Fill in F:I results for rows 134-148 (new IOU2 / Thresh2 columns of results)
$data = @(
    , @(134, 0.308, 0.96, 0.366, 0.92)
    , @(135, 0.317, 0.92, 0.376, 0.92)
    , @(136, 0.307, 0.79, 0.364, 0.79)
    , @(137, 0.309, 0.83, 0.364, 0.79)
    , @(138, 0.321, 0.87, 0.368, 0.83)
    , @(139, 0.323, 0.79, 0.378, 0.79)
    , @(140, 0.334, 0.83, 0.387, 0.83)
    , @(141, 0.321, 0.79, 0.374, 0.79)
    , @(142, 0.315, 0.66, 0.369, 0.66)
    , @(143, 0.303, 0.75, 0.368, 0.75)
    , @(144, 0.015, 0.13, 0.03, 0.13)
    , @(145, 0.015, 0.13, 0.03, 0.13)
    , @(146, 0.314, 0.67, 0.373, 0.67)
    , @(147, 0.105, 0.21, 0.16, 0.21)
    , @(148, 0.305, 0.54, 0.355, 0.54)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 6).Value = $entry[1]
    $ws.Cells.Item($row, 7).Value = $entry[2]
    $ws.Cells.Item($row, 8).Value = $entry[3]
    $ws.Cells.Item($row, 9).Value = $entry[4]
}

# Append a new results row (168) to the last experiment block
$ws.Range("A168").Value = 0.00004
$ws.Range("B168").Value = 0.0004
$ws.Range("C168").Value = "run118"

# Update view state: active selection cell
$ws.Range("B168").Select() | Out-Null
